# Updated cryptos list with GitHub Actions - refresh market data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference a never-edited, default-styled cell so numeric-looking
# price strings we re-enter as text do not pick up a stray quote-prefix style
$plainStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "95.020.85"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "3.579.81"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("D5").Value = "'2.25"
$ws.Range("E5").Value = "  +17.34%  "
$ws.Range("D6").Value = "'223.68"
$ws.Range("E6").Value = "  -5.46%  "
$ws.Range("D7").Value = "'631.27"
$ws.Range("E7").Value = "  -3.90%  "
$ws.Range("E8").Value = "  -4.18%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "3.577.51"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("D12").Value = "'45.44"
$ws.Range("E12").Value = "  +2.87%  "
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("E14").Value = "  -5.09%  "
$ws.Range("E15").Value = "  -5.56%  "
$ws.Range("D16").Value = "4.248.02"
$ws.Range("E16").Value = "  -2.93%  "
$ws.Range("D17").Value = "94.696.99"
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").Value = "'8.70"
$ws.Range("E18").Value = "  -4.72%  "
$ws.Range("D19").Value = "'19.69"
$ws.Range("E19").Value = "  +5.15%  "
$ws.Range("D20").Value = "3.572.32"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").Value = "'12.75"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").Value = "'0.501"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").Value = "'497.22"
$ws.Range("E23").Value = "  -4.17%  "
$ws.Range("D24").Value = "'3.20"
$ws.Range("E24").Value = "  -6.10%  "
$ws.Range("D25").Value = "'0.234"
$ws.Range("E25").Value = "  +13.92%  "
$ws.Range("D26").Value = "'115.80"
$ws.Range("E26").Value = "  +14.30%  "
$ws.Range("E27").Value = "  -4.79%  "
$ws.Range("E28").Value = "  -3.74%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'12.38"
$ws.Range("E29").Value = "  -7.17%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'12.55"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.86"
$ws.Range("E31").Value = "  -5.13%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D34").Value = "'0.177"
$ws.Range("E34").Value = "  -6.30%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.73"
$ws.Range("E35").Value = "  -7.00%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'31.27"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.574"
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("B38").Value = "USDe"
$ws.Range("C38").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'582.66"
$ws.Range("E39").Value = "  -9.89%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "'8.18"
$ws.Range("E40").Value = "  -7.25%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'6.70"
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'40.31"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.157"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.460"
$ws.Range("E44").Value = "  -6.24%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "'1.88"
$ws.Range("E45").Value = "  -9.04%  "
$ws.Range("D46").Value = "'0.0461"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'0.906"
$ws.Range("E47").Value = "  -5.43%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'23.41"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("B49").Value = "MantraDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D49").Value = "'3.59"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'217.35"
$ws.Range("E50").Value = "  +6.58%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'8.38"
$ws.Range("E51").Value = "  -2.98%  "

# Restore the default (non quote-prefixed) style on those cells
$ws.Range("D5").Style = $plainStyle
$ws.Range("D6").Style = $plainStyle
$ws.Range("D7").Style = $plainStyle
$ws.Range("D12").Style = $plainStyle
$ws.Range("D18").Style = $plainStyle
$ws.Range("D19").Style = $plainStyle
$ws.Range("D21").Style = $plainStyle
$ws.Range("D22").Style = $plainStyle
$ws.Range("D23").Style = $plainStyle
$ws.Range("D24").Style = $plainStyle
$ws.Range("D25").Style = $plainStyle
$ws.Range("D26").Style = $plainStyle
$ws.Range("D29").Style = $plainStyle
$ws.Range("D30").Style = $plainStyle
$ws.Range("D31").Style = $plainStyle
$ws.Range("D32").Style = $plainStyle
$ws.Range("D34").Style = $plainStyle
$ws.Range("D35").Style = $plainStyle
$ws.Range("D36").Style = $plainStyle
$ws.Range("D37").Style = $plainStyle
$ws.Range("D38").Style = $plainStyle
$ws.Range("D39").Style = $plainStyle
$ws.Range("D40").Style = $plainStyle
$ws.Range("D41").Style = $plainStyle
$ws.Range("D42").Style = $plainStyle
$ws.Range("D43").Style = $plainStyle
$ws.Range("D44").Style = $plainStyle
$ws.Range("D45").Style = $plainStyle
$ws.Range("D46").Style = $plainStyle
$ws.Range("D47").Style = $plainStyle
$ws.Range("D48").Style = $plainStyle
$ws.Range("D49").Style = $plainStyle
$ws.Range("D50").Style = $plainStyle
$ws.Range("D51").Style = $plainStyle
